$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column K (최종점수) values - decreased by 0.3
$ws.Range("K2").Value = 59.2
$ws.Range("K3").Value = 55.2
$ws.Range("K4").Value = 53.2
$ws.Range("K5").Value = 52.2

# Update column N (MACRO_SCORE) values - new value for all rows
$ws.Range("N2").Value = 50.60178744571824
$ws.Range("N3").Value = 50.60178744571824
$ws.Range("N4").Value = 50.60178744571824
$ws.Range("N5").Value = 50.60178744571824
